# "On and Off Switch implemented"
# Adds a new "Status" column (I) to the tracker sheet with an On/Off
# value per row, matching whichever row currently carries the "active"
# color tag (PALE_GREEN/PALE_BLUE rows -> On/Off accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell I1: "Status", bold header style like the rest of row 1 ---
$ws.Range("I1").Value = "Status"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").Font.ThemeFont = 1

# --- New data column I2:I6, plain (unstyled) cells like the rest of the data rows ---
$ws.Range("I2").Value = "On"
$ws.Range("I3").Value = "Off"
$ws.Range("I4").Value = "Off"
$ws.Range("I5").Value = "On"
$ws.Range("I6").Value = "On"

# Scroll the sheet view so column F is the leftmost visible column.
$ws.Range("F1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
